$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Decrement the "Sprint #" value (column C) by 1 for every backlog row
# (rows 3-29), except row 19 which is left unchanged.
foreach ($r in 3..29) {
    if ($r -eq 19) { continue }
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    $cell.Value = $current - 1
}

# Update the active selection to B4
$ws.Range("B4").Select()
